# Adjust Enemy Attack Delay
# - Halve the ATTACKDELAY values for several enemies on the "Enemy" sheet.
# - Add a new ACCURANCY value ("100;0;100;0") to several items on the "Item" sheet.
# - Misc view-state changes (active sheet/tab, selections, column widths) left
#   over from the author's interactive editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Item" sheet: populate column H (ACCURANCY) for a batch of items that didn't
# have an accuracy value yet.
# ---------------------------------------------------------------------------
$itemSheet = $wb.Worksheets.Item("Item")

$accuracyValue = "100;0;100;0"
foreach ($r in 7,9,10,11,12,15,16,17,19,20) {
    $itemSheet.Cells.Item($r, 8).Value = $accuracyValue
}

# widen column G ("AMMO") a bit so the longer values introduced above are readable
$itemSheet.Columns.Item(7).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# "Enemy" sheet: cut ATTACKDELAY (column E) roughly in half for most enemies.
# ---------------------------------------------------------------------------
$enemySheet = $wb.Worksheets.Item("Enemy")

$enemySheet.Range("E2").Value = 1      # SWORD        2   -> 1
$enemySheet.Range("E3").Value = 0.5    # GUN          1   -> 0.5
$enemySheet.Range("E4").Value = 1.5    # CAR          3   -> 1.5
$enemySheet.Range("E5").Value = 1      # PARATROOPER  4   -> 1
$enemySheet.Range("E7").Value = 0.1    # FLAMETHROWER 0.2 -> 0.1
$enemySheet.Range("E8").Value = 7.5    # TANK         15  -> 7.5
$enemySheet.Range("E9").Value = 2.5    # HELICOPTER   10  -> 2.5
$enemySheet.Range("E10").Value = 1     # ROBOT        2   -> 1

# widen column E ("ATTACKDELAY") to fit the new decimal values
$enemySheet.Columns.Item(5).ColumnWidth = 15.15

# ---------------------------------------------------------------------------
# View-state: selections on both sheets, and make "Enemy" the active tab.
# ---------------------------------------------------------------------------
$itemSheet.Range("H14").Select() | Out-Null

$enemySheet.Activate() | Out-Null
$enemySheet.Range("G8").Select() | Out-Null
